$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before the old row 101 (old rows 101-157 shift down to 103-159)
$ws.Rows.Item(101).Insert()
$ws.Rows.Item(102).Insert()

# New row 101: Murcott / Primera, Terminal Hortofruticola Agro Chillan, Ñuble, 2021-10-20
$ws.Cells.Item(101, 1).Value = 7
$ws.Cells.Item(101, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(101, 3).Value = "Ñuble"
$ws.Cells.Item(101, 4).Value = (Get-Date -Year 2021 -Month 10 -Day 20 -Hour 0 -Minute 0 -Second 0)
$ws.Cells.Item(101, 5).Value = 16
$ws.Cells.Item(101, 6).Value = "Fruta"
$ws.Cells.Item(101, 7).Value = 100102
$ws.Cells.Item(101, 8).Value = "Cítricos"
$ws.Cells.Item(101, 9).Value = 100102004
$ws.Cells.Item(101, 10).Value = "Mandarina"
$ws.Cells.Item(101, 11).Value = "Murcott"
$ws.Cells.Item(101, 12).Value = "Primera"
$ws.Cells.Item(101, 13).Value = 160
$ws.Cells.Item(101, 14).Value = 7000
$ws.Cells.Item(101, 15).Value = 7500
$ws.Cells.Item(101, 16).Value = 7250
$ws.Cells.Item(101, 17).Value = "$/caja 18 kilos"
$ws.Cells.Item(101, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(101, 19).Value = 403
$ws.Cells.Item(101, 20).Value = 18

# New row 102: Murcott / Segunda, Terminal Hortofruticola Agro Chillan, Ñuble, 2021-10-20
$ws.Cells.Item(102, 1).Value = 7
$ws.Cells.Item(102, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(102, 3).Value = "Ñuble"
$ws.Cells.Item(102, 4).Value = (Get-Date -Year 2021 -Month 10 -Day 20 -Hour 0 -Minute 0 -Second 0)
$ws.Cells.Item(102, 5).Value = 16
$ws.Cells.Item(102, 6).Value = "Fruta"
$ws.Cells.Item(102, 7).Value = 100102
$ws.Cells.Item(102, 8).Value = "Cítricos"
$ws.Cells.Item(102, 9).Value = 100102004
$ws.Cells.Item(102, 10).Value = "Mandarina"
$ws.Cells.Item(102, 11).Value = "Murcott"
$ws.Cells.Item(102, 12).Value = "Segunda"
$ws.Cells.Item(102, 13).Value = 160
$ws.Cells.Item(102, 14).Value = 6000
$ws.Cells.Item(102, 15).Value = 6500
$ws.Cells.Item(102, 16).Value = 6250
$ws.Cells.Item(102, 17).Value = "$/caja 18 kilos"
$ws.Cells.Item(102, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(102, 19).Value = 347
$ws.Cells.Item(102, 20).Value = 18

Write-Host "Done inserting rows 101-102"
